$d = $word.ActiveDocument

# Word inherits the *direct* run formatting of whichever paragraph a new
# paragraph is split off from (InsertParagraphBefore/After). The title
# paragraph ("Dheeraj Chand") carries direct Bold + 28pt formatting, so
# building the new contact-info paragraph right next to it would leak that
# formatting into the new run. Instead, build the new paragraph next to a
# paragraph that carries no direct character formatting (the plain summary
# paragraph), then relocate (cut/paste) the whole paragraph - including its
# paragraph mark - to the correct spot right after the title. Cut/paste moves
# the paragraph as-is, so no stray pStyle/rPr gets introduced.

$donor = $d.Paragraphs(3)
$donor.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs(4)
$newPara.Range.Text = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"
$newPara.Format.Alignment = 1

$cutRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$cutRange.Cut()

$titlePara = $d.Paragraphs(1)
$insertPoint = $titlePara.Range.End
$target = $d.Range($insertPoint, $insertPoint)
$target.Paste()
